$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.868.76'
$ws.Range("E2").Value = '  +0.92%  '

$ws.Range("D3").Value = '2.220.03'
$ws.Range("E3").Value = '  -0.27%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '292.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.65%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '87.31'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.66%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.516'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("E9").Value = '  -0.04%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '30.47'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.03%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0785'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.93%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.54'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.94%  '

$ws.Range("D15").Value = '2.561.93'
$ws.Range("E15").Value = '  -0.34%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.08'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.60%  '

$ws.Range("D17").Value = '2.224.60'
$ws.Range("E17").Value = '  -0.03%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.731'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.60%  '

$ws.Range("D19").Value = '39.804.70'
$ws.Range("E19").Value = '  +0.98%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.50'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +11.23%  '

$ws.Range("E21").Value = '  +0.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.84'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.80'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.70'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.81%  '

$ws.Range("E25").Value = '  -0.03%  '

$ws.Range("E26").Value = '  +2.29%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.84'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.74%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.81'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.24%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.20'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.17%  '

$ws.Range("E30").Value = '  +1.24%  '

$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '152.66'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.57%  '

$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.73'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.36%  '

$ws.Range("E33").Value = '  -0.01%  '

$ws.Range("E34").Value = '  +1.93%  '

$ws.Range("E35").Value = '  +2.57%  '

$ws.Range("E37").Value = '  +5.96%  '

$ws.Range("E38").Value = '  +1.17%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.92'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.79%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0991'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.19%  '

$ws.Range("E41").Value = '  +2.21%  '

$ws.Range("D42").Value = '2.100.52'
$ws.Range("E42").Value = '  +9.09%  '

$ws.Range("E43").Value = '  +2.81%  '

$ws.Range("E44").Value = '  +6.41%  '

$ws.Range("E45").Value = '  +2.83%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.01'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.65%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.71'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.51%  '

$ws.Range("E48").Value = '  +1.25%  '

$ws.Range("D49").Value = '2.436.22'
$ws.Range("E49").Value = '  -0.12%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '70.88'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.91%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '89.29'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.96%  '
